$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.285.98'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.919.28'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.63%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7442'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.97%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.75%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.22%  '

$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.31'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.71%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3134'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.61%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06964'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07988'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.41%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7722'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.96%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.916.98'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.310'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.43%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.63'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.20%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.266.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.22%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.22'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '246.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.71%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.834'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007871'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.36%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.34%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.168.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.77%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.21%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.627'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.92%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.414'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.54%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.68'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.46%  '

$ws.Range('E28').Value = '  -2.80%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.151'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.80%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.359'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.543'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.348'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.59%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.075'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.81%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05185'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.30%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.302'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.08%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7484'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.25%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.774'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01934'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.16%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.790'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.50%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.419'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.25%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.84%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4469'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.70%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.943'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.28%  '

$ws.Range('E44').Value = '  +0.18%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8383'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.89%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.653'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.64%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.36'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.05%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.859'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.92%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.081.24'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.35%  '

$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.22%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1221'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.09%  '
